$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '30.536.63'
Set-TextValue "E2" '  +0.13%  '

Set-TextValue "D3" '1.918.51'
Set-TextValue "E3" '  -0.21%  '

Set-TextValue "D4" '1.001'
Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '245.29'
Set-TextValue "E5" '  +0.94%  '

Set-TextValue "D6" '1.001'
Set-TextValue "E6" '  -0.06%  '

Set-TextValue "D7" '0.4794'
Set-TextValue "E7" '  +2.06%  '

Set-TextValue "D8" '0.2888'
Set-TextValue "E8" '  +0.72%  '

Set-TextValue "D9" '0.06725'
Set-TextValue "E9" '  -0.13%  '

Set-TextValue "D10" '110.02'

Set-TextValue "D11" '19.11'
Set-TextValue "E11" '  +4.51%  '

Set-TextValue "D12" '1.916.30'
Set-TextValue "E12" '  -0.26%  '

Set-TextValue "D13" '0.07573'
Set-TextValue "E13" '  -2.33%  '

Set-TextValue "D14" '5.254'
Set-TextValue "E14" '  -0.90%  '

Set-TextValue "D15" '0.6664'
Set-TextValue "E15" '  +1.19%  '

Set-TextValue "D16" '298.34'
Set-TextValue "E16" '  +2.79%  '

Set-TextValue "D17" '30.505.83'
Set-TextValue "E17" '  +0.05%  '

Set-TextValue "D18" '13.02'
Set-TextValue "E18" '  +0.65%  '

Set-TextValue "B19" 'Uniswap'
Set-TextValue "C19" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D19" '5.582'
Set-TextValue "E19" '  +6.18%  '

Set-TextValue "B20" 'Dai'
Set-TextValue "C20" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D20" '1.000'
Set-TextValue "E20" '  -0.09%  '

Set-TextValue "B21" 'ShibaInu'
Set-TextValue "C21" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D21" '0.000007572'
Set-TextValue "E21" '  -0.17%  '

Set-TextValue "D22" '2.158.45'
Set-TextValue "E22" '  +0.39%  '

Set-TextValue "D23" '1.001'
Set-TextValue "E23" '  +0.06%  '

Set-TextValue "D24" '6.417'
Set-TextValue "E24" '  +3.53%  '

Set-TextValue "D25" '9.470'
Set-TextValue "E25" '  +0.96%  '

Set-TextValue "D26" '164.47'
Set-TextValue "E26" '  -2.81%  '

Set-TextValue "D27" '20.31'
Set-TextValue "E27" '  -4.52%  '

Set-TextValue "D28" '2.117'
Set-TextValue "E28" '  -0.43%  '

Set-TextValue "D29" '0.1079'
Set-TextValue "E29" '  +1.12%  '

Set-TextValue "D30" '1.398'
Set-TextValue "E30" '  +2.13%  '

Set-TextValue "D31" '4.163'
Set-TextValue "E31" '  -0.26%  '

Set-TextValue "D32" '4.022'
Set-TextValue "E32" '  +1.08%  '

Set-TextValue "D33" '0.04994'
Set-TextValue "E33" '  -0.71%  '

Set-TextValue "D34" '0.7352'
Set-TextValue "E34" '  -0.63%  '

Set-TextValue "E35" '  -1.15%  '

Set-TextValue "B36" 'Frax'
Set-TextValue "C36" 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue "D36" '0.9990'
Set-TextValue "E36" '  -0.04%  '

Set-TextValue "B37" 'VeChain'
Set-TextValue "C37" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D37" '0.02045'
Set-TextValue "E37" '  -2.14%  '

Set-TextValue "D38" '2.724'
Set-TextValue "E38" '  +0.22%  '

Set-TextValue "D39" '2.682'
Set-TextValue "E39" '  +0.11%  '

Set-TextValue "D40" '110.61'
Set-TextValue "E40" '  +0.39%  '

Set-TextValue "D41" '2.022'
Set-TextValue "E41" '  -1.80%  '

Set-TextValue "D42" '0.4430'
Set-TextValue "E42" '  +4.26%  '

Set-TextValue "D43" '0.8658'
Set-TextValue "E43" '  -0.67%  '

Set-TextValue "B44" 'Aave'
Set-TextValue "C44" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D44" '71.77'
Set-TextValue "E44" '  +6.77%  '

Set-TextValue "B45" 'FraxShare'
Set-TextValue "C45" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D45" '5.907'
Set-TextValue "E45" '  +0.82%  '

Set-TextValue "D46" '1.001'
Set-TextValue "E46" '  -0.07%  '

Set-TextValue "D47" '49.53'
Set-TextValue "E47" '  -0.30%  '

Set-TextValue "D48" '7.280'
Set-TextValue "E48" '  +1.18%  '

Set-TextValue "D49" '9.315'
Set-TextValue "E49" '  +0.87%  '

Set-TextValue "D50" '0.1232'
Set-TextValue "E50" '  +1.39%  '

Set-TextValue "D51" '0.2539'
Set-TextValue "E51" '  +3.16%  '
